$d = $word.ActiveDocument

# 1. Date change in header paragraph
$d.Paragraphs.Item(1).Range.Text = "⚡️🚀המאמר היומי של מייק 24.06.24:⚡️🚀"

# 2. Title
$d.Paragraphs.Item(2).Range.Text = "Are you still on track!? Catching LLM Task Drift with Activations"

# 3. Intro paragraph
$d.Paragraphs.Item(3).Range.Text = "הסקירה הזו הולכת להיות קצרה כי הרעיון העיקרי של המאמר הוא די פשוט ואינטואיטיבי. אתם מדברים עם מודלי שפה שלכם באמצעות שאילתות שבד״כ נקראות פרומפטים שהמודל עונה לכם. אבל מה קורה אם מודל השפה שלכם מחובר לעוד כלי שמגנרט בשבילו פרומפטים למשל בהתבסס על תוצאה של איזשהו חישוב על הפלט של מודל אחר או מתבסס על RAG או אולי אפילו תלוי בתוצאות חיפוש באינטרנט."

# 4. Problem paragraph
$d.Paragraphs.Item(4).Range.Text = "כמובן שגנרוט אוטומטי של פרומפט יכול להתפקשש (באגים, אולי פעילות זדונית) ואז יחד עם שאלה לגיטימית המודל מקבל תופסת לא קשורה. בעיה ידועה, אה?"

# 5. Approach paragraph
$d.Paragraphs.Item(5).Range.Text = "אז המאמר שבנידון חקר את האקטיבציות של שכבות המודל (טרנספורמר כמובן) ומצאו הבדלים משמעותיים בין האקטיבציות הנוצרות על ידי שאלה לגיטימית לבין אלו שנוצרו עם שאלה ״מקושקשת״. ואז הם בנו דאטהסט של שאלות טובות ושאלות מורעלות ואימנו מודל (קטן) שיודע להבדיל בין האקטיבציות של שאלות הטובות והלא טובות. המחברים לוקחים אקטיבציות של הטוקן האחרון של הפרומפט (השאלה) המלא"

# 6. Methods paragraph
$d.Paragraphs.Item(6).Range.Text = "הם ניסו שתי שיטות: אחת היא אימון של שכבה לינארית המפרידה בין ייצוגים טובים ומורעלים. השיטה השניה שהם מנסים נקראת metric learning שבמילים פשוטות מנסה ללמוד ייצוג (המופק על ידי המודל ״המבדיל״) המקרב ייצוגים של העוגן (התחלת השאלה) עם השאלה הטובה ומרחיק אותו מהייצוג של השאלה המורעלת (התוספת המורעלת). אם מצליחים ב-metric learning אז בקלות אפשר לתפור שכבה לינארית המבדילה בין הטובים ללא טובים."

# 9. Link paragraph (update before removing the two paragraphs in between,
# so indices 7/8/9 remain valid)
$d.Paragraphs.Item(9).Range.Text = "https://arxiv.org/pdf/2406.00799"

# 7 & 8. Remove the LoRa paragraph and the averaging paragraph entirely
$d.Paragraphs.Item(8).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()
